# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (holdings detail) right before the
#    "总计" (summary) sheet.
# 2. Add a new top row to the "总计" sheet summarizing the 2022-Q1 quarter,
#    shifting the previously-existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: add the "2022-Q1" worksheet just before "总计"
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Fund codes / names / the decimal-looking figures are all stored as *text*
# in the source data (leading zeros in codes, trailing zeros in numbers),
# so format columns B:G as Text before writing any values into them.
$newSheet.Range("B1:G10").NumberFormat = "@"

# Borrow the existing formatting used on other quarter sheets: the bold
# centered/bordered style used for the header row (B1:H1) and for the
# row-index column (A2:A10).
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Range("A2").Copy()
$newSheet.Range("A2:A10").PasteSpecial(-4122)   # xlPasteFormats
$q4.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)    # xlPasteFormats

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows
$data = @(
    @(0, "012588", "南方港股通优势企业混合型证券投资基金A", "37.54", "71.00", "2.00", "0.7508", 9),
    @(1, "010761", "华商甄选回报混合", "20.63", "93.93", "2.69", "0.5549", 9),
    @(2, "013626", "华夏周期驱动混合A", "1.96", "90.84", "2.34", "0.0459", 10),
    @(3, "012589", "南方港股通优势企业混合型证券投资基金C", "2.05", "71.00", "2.00", "0.0410", 9),
    @(4, "013627", "华夏周期驱动混合C", "1.18", "90.84", "2.34", "0.0276", 10),
    @(5, "004098", "前海开源港股通股息率50强股票", "0.34", "88.92", "3.78", "0.0129", 4),
    @(6, "501303", "广发港股通恒生综合中型股指数(LOF)A", "0.34", "92.39", "1.47", "0.0050", 7),
    @(7, "004996", "广发港股通恒生综合中型股指数(LOF)C", "0.11", "92.39", "1.47", "0.0016", 7),
    @(8, "160922", "大成恒生综合中小型股指数(QDII-LOF)A", "0.10", "92.44", "1.17", "0.0012", 6)
)

$rowIdx = 2
foreach ($rec in $data) {
    $newSheet.Cells.Item($rowIdx, 1).Value = $rec[0]
    $newSheet.Cells.Item($rowIdx, 2).Value = $rec[1]
    $newSheet.Cells.Item($rowIdx, 3).Value = $rec[2]
    $newSheet.Cells.Item($rowIdx, 4).Value = $rec[3]
    $newSheet.Cells.Item($rowIdx, 5).Value = $rec[4]
    $newSheet.Cells.Item($rowIdx, 6).Value = $rec[5]
    $newSheet.Cells.Item($rowIdx, 7).Value = $rec[6]
    $newSheet.Cells.Item($rowIdx, 8).Value = $rec[7]
    $rowIdx = $rowIdx + 1
}

# ---------------------------------------------------------------------------
# Step 2: add the 2022-Q1 summary row to "总计", pushing old rows down
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")

# Read the current values (bottom rows first isn't required since we grab
# everything up-front before writing anything back). The "A" column is just
# a 0-based row index, so only B/C/D (the actual data) need to be preserved.
$b2 = $ws.Range("B2").Value2
$c2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2

$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2
$d3 = $ws.Range("D3").Value2

$b4 = $ws.Range("B4").Value2
$c4 = $ws.Range("C4").Value2
$d4 = $ws.Range("D4").Value2

$b5 = $ws.Range("B5").Value2
$c5 = $ws.Range("C5").Value2
$d5 = $ws.Range("D5").Value2

# Give the brand-new row 6 the same formatting as the existing index column.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)   # xlPasteFormats

# Shift rows 2-5 down into rows 3-6, renumbering the index column 1..4.
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = $b5
$ws.Range("C6").Value = $c5
$ws.Range("D6").Value = $d5

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = $b4
$ws.Range("C5").Value = $c4
$ws.Range("D5").Value = $d4

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = $b3
$ws.Range("C4").Value = $c3
$ws.Range("D4").Value = $d3

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = $b2
$ws.Range("C3").Value = $c2
$ws.Range("D3").Value = $d2

# Insert the new 2022-Q1 totals into row 2.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 1.44
